$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Map of (row, col) -> new value, using 1-based Word COM indices.
# The data rows in this 20x5 table are rows 1, 5, 10, 15, 20.
$updates = @{
    "1,1"  = "435×3=1305"
    "1,2"  = "992×4=3968"
    "1,3"  = "839×4=3356"
    "1,4"  = "916×2=1832"
    "1,5"  = "329×2=658"

    "5,1"  = "157×4=628"
    "5,2"  = "637×4=2548"
    "5,3"  = "873×8=6984"
    "5,4"  = "423×7=2961"
    "5,5"  = "173×6=1038"

    "10,1" = "261×2=522"
    "10,2" = "774×3=2322"
    "10,3" = "403×9=3627"
    "10,4" = "630×7=4410"
    "10,5" = "298×7=2086"

    "15,1" = "997×7=6979"
    "15,2" = "199×3=597"
    "15,3" = "821×2=1642"
    "15,4" = "276×6=1656"
    "15,5" = "414×9=3726"

    "20,1" = "528×6=3168"
    "20,2" = "361×6=2166"
    "20,3" = "412×2=824"
    "20,4" = "911×7=6377"
    "20,5" = "609×5=3045"
}

foreach ($key in $updates.Keys) {
    $parts = $key.Split(",")
    $row = [int]$parts[0]
    $col = [int]$parts[1]
    $cell = $t.Cell($row, $col)
    $cell.Range.Text = $updates[$key]
}
